{"js": "// Update the worksheet date and the 25 multiplication problems/answers\n// (three-digit number multiplied by one-digit number) to the new values.\nconst replacements = [\n  [\"2025-12-06 Saturday\", \"2025-12-07 Sunday\"],\n  [\"201\u00d76=1206\", \"215\u00d78=1720\"],\n  [\"937\u00d76=5622\", \"449\u00d72=898\"],\n  [\"274\u00d77=1918\", \"120\u00d76=720\"],\n  [\"250\u00d78=2000\", \"963\u00d78=7704\"],\n  [\"947\u00d78=7576\", \"275\u00d77=1925\"],\n  [\"718\u00d78=5744\", \"252\u00d74=1008\"],\n  [\"448\u00d76=2688\", \"559\u00d78=4472\"],\n  [\"625\u00d78=5000\", \"750\u00d74=3000\"],\n  [\"713\u00d72=1426\", \"114\u00d75=570\"],\n  [\"973\u00d75=4865\", \"454\u00d74=1816\"],\n  [\"113\u00d77=791\", \"118\u00d72=236\"],\n  [\"984\u00d76=5904\", \"664\u00d72=1328\"],\n  [\"349\u00d75=1745\", \"685\u00d76=4110\"],\n  [\"535\u00d72=1070\", \"339\u00d76=2034\"],\n  [\"803\u00d75=4015\", \"354\u00d74=1416\"],\n  [\"360\u00d79=3240\", \"664\u00d79=5976\"],\n  [\"589\u00d72=1178\", \"863\u00d77=6041\"],\n  [\"647\u00d78=5176\", \"446\u00d76=2676\"],\n  [\"920\u00d73=2760\", \"880\u00d75=4400\"],\n  [\"633\u00d75=3165\", \"342\u00d72=684\"],\n  [\"291\u00d79=2619\", \"949\u00d77=6643\"],\n  [\"418\u00d73=1254\", \"659\u00d72=1318\"],\n  [\"624\u00d73=1872\", \"831\u00d75=4155\"],\n  [\"214\u00d73=642\", \"485\u00d73=1455\"],\n  [\"786\u00d79=7074\", \"638\u00d72=1276\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 multiplication problems/answers\n# (three-digit number multiplied by one-digit number) to the new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-06 Saturday\", \"2025-12-07 Sunday\"),\n    @(\"201\u00d76=1206\", \"215\u00d78=1720\"),\n    @(\"937\u00d76=5622\", \"449\u00d72=898\"),\n    @(\"274\u00d77=1918\", \"120\u00d76=720\"),\n    @(\"250\u00d78=2000\", \"963\u00d78=7704\"),\n    @(\"947\u00d78=7576\", \"275\u00d77=1925\"),\n    @(\"718\u00d78=5744\", \"252\u00d74=1008\"),\n    @(\"448\u00d76=2688\", \"559\u00d78=4472\"),\n    @(\"625\u00d78=5000\", \"750\u00d74=3000\"),\n    @(\"713\u00d72=1426\", \"114\u00d75=570\"),\n    @(\"973\u00d75=4865\", \"454\u00d74=1816\"),\n    @(\"113\u00d77=791\", \"118\u00d72=236\"),\n    @(\"984\u00d76=5904\", \"664\u00d72=1328\"),\n    @(\"349\u00d75=1745\", \"685\u00d76=4110\"),\n    @(\"535\u00d72=1070\", \"339\u00d76=2034\"),\n    @(\"803\u00d75=4015\", \"354\u00d74=1416\"),\n    @(\"360\u00d79=3240\", \"664\u00d79=5976\"),\n    @(\"589\u00d72=1178\", \"863\u00d77=6041\"),\n    @(\"647\u00d78=5176\", \"446\u00d76=2676\"),\n    @(\"920\u00d73=2760\", \"880\u00d75=4400\"),\n    @(\"633\u00d75=3165\", \"342\u00d72=684\"),\n    @(\"291\u00d79=2619\", \"949\u00d77=6643\"),\n    @(\"418\u00d73=1254\", \"659\u00d72=1318\"),\n    @(\"624\u00d73=1872\", \"831\u00d75=4155\"),\n    @(\"214\u00d73=642\", \"485\u00d73=1455\"),\n    @(\"786\u00d79=7074\", \"638\u00d72=1276\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
